$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D2").Value = "'57.685.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "'3.015.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'512.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'140.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.438"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").Value = "'7.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "'0.111"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").Value = "'0.367"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("D12").Value = "'3.527.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "'26.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.80%  "
$ws.Range("D15").Value = "'0.0000165"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.90%  "
$ws.Range("D16").Value = "'57.616.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").Value = "'6.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.07%  "
$ws.Range("D18").Value = "'3.014.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'12.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.30%  "
$ws.Range("D20").Value = "'8.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").Value = "'332.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'0.500"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.85%  "
$ws.Range("D24").Value = "'64.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "'0.0₃0929"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.90%  "
$ws.Range("D28").Value = "'6.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.09%  "
$ws.Range("D29").Value = "'7.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.37%  "
$ws.Range("D30").Value = "'1.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.02%  "
$ws.Range("E31").Value = "  -5.11%  "
$ws.Range("D32").Value = "'20.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").Value = "'4.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.02%  "
$ws.Range("D34").Value = "'154.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").Value = "'5.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.28%  "
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").Value = "'24.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").Value = "'0.0679"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").Value = "'3.048.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "'37.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").Value = "'3.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.70%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "'2.242.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "'1.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("D46").Value = "'0.989"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").Value = "'6.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.14%  "
$ws.Range("D48").Value = "'0.0241"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").Value = "'19.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("D50").Value = "'1.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.55%  "
$ws.Range("D51").Value = "'0.0895"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.97%  "
